# Apply "repull data, push all data, mean calculation" changes
# Updates to column F (dSF) values on Sheet1, reflecting repulled data.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("F2").Value = -4
$ws.Range("F3").Value = 10
$ws.Range("F4").Value = -2
$ws.Range("F5").Value = -8
$ws.Range("F9").Value = 13
$ws.Range("F11").Value = 0
$ws.Range("F13").Value = -7
